$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 4 ("Git commands (console)") - "git pull" bullet:
# split the trailing sentence into two runs and add "and combine it with
# our local work" to the end.
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$tr4 = $shp4.TextFrame.TextRange
$para4 = $tr4.Paragraphs(5)
$run4 = $para4.Runs(2)
$run4.Text = "- get all new work from the remote central repository that is not in our local repository and combine it "
[void]$run4.InsertAfter("with our local work")

# ---------------------------------------------------------------------------
# Slide 8 ("git pull") - "Content Placeholder 2":
# grow the placeholder, turn on shrink-text-on-overflow autofit and extend
# the bullet text.
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(2)
$shp8.Width = 5097174 / 12700
$shp8.Height = 2054881 / 12700
$shp8.TextFrame.AutoSize = 2

$tr8 = $shp8.TextFrame.TextRange
$para8 = $tr8.Paragraphs(1)
$run8 = $para8.Runs(2)
$run8.Text = "- get all new work from the remote central repository that is not in our local repository and combine with our local work."
